$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (user 5163876201 / VaLaK_DEMON)
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "2026-02-12T10:37:52.409468+00:00"
$ws.Range("E5").Value = "I'm looking for intern"
$ws.Range("F5").Value = "yes"

# Row 6 (user 1900918712 / lunaticbeast12)
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "2026-02-12T10:38:30.743047+00:00"
$ws.Range("E6").Value = "ys bro"
$ws.Range("F6").Value = "yes"
$ws.Range("H6").Value = $true
$ws.Range("J6").Value = "success"
